$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text update: F1 "total_2021" -> "total_2022" ---
$ws.Range("F1").Value = "total_2022"

# --- Data rows 2-14: updated figures (edition 10 of the PDF) ---
$ws.Range("B2").Value = 55404
$ws.Range("C2").Value = 55146
$ws.Range("D2").Value = 414063
$ws.Range("E2").Value = 14966
$ws.Range("F2").Value = 539579
$ws.Range("B3").Value = 50790
$ws.Range("C3").Value = 50617
$ws.Range("D3").Value = 408170
$ws.Range("E3").Value = 12509
$ws.Range("F3").Value = 522086
$ws.Range("B4").Value = 59201
$ws.Range("C4").Value = 59017
$ws.Range("D4").Value = 515945
$ws.Range("E4").Value = 12293
$ws.Range("F4").Value = 646456
$ws.Range("B5").Value = 63347
$ws.Range("C5").Value = 63158
$ws.Range("D5").Value = 604074
$ws.Range("E5").Value = 11410
$ws.Range("F5").Value = 741989
$ws.Range("B6").Value = 71101
$ws.Range("C6").Value = 70824
$ws.Range("D6").Value = 695949
$ws.Range("E6").Value = 13273
$ws.Range("F6").Value = 851147
$ws.Range("B7").Value = 77150
$ws.Range("C7").Value = 76794
$ws.Range("D7").Value = 725748
$ws.Range("E7").Value = 13977
$ws.Range("F7").Value = 893669
$ws.Range("B8").Value = 87541
$ws.Range("C8").Value = 87643
$ws.Range("D8").Value = 759745
$ws.Range("E8").Value = 15615
$ws.Range("F8").Value = 950544
$ws.Range("B9").Value = 88379
$ws.Range("C9").Value = 88363
$ws.Range("D9").Value = 755485
$ws.Range("E9").Value = 15897
$ws.Range("F9").Value = 948124
$ws.Range("B10").Value = 81577
$ws.Range("C10").Value = 81483
$ws.Range("D10").Value = 721100
$ws.Range("E10").Value = 15117
$ws.Range("F10").Value = 899277
$ws.Range("B11").Value = 79570
$ws.Range("C11").Value = 80012
$ws.Range("D11").Value = 684229
$ws.Range("E11").Value = 15033
$ws.Range("F11").Value = 858844
$ws.Range("B12").Value = 70868
$ws.Range("C12").Value = 70939
$ws.Range("D12").Value = 536747
$ws.Range("E12").Value = 14772
$ws.Range("F12").Value = 693326
$ws.Range("B13").Value = 74021
$ws.Range("C13").Value = 74221
$ws.Range("D13").Value = 527941
$ws.Range("E13").Value = 16469
$ws.Range("F13").Value = 692652
$ws.Range("B14").Value = 858949
$ws.Range("C14").Value = 858217
$ws.Range("D14").Value = 7349196
$ws.Range("E14").Value = 171331
$ws.Range("F14").Value = 9237693


# --- Recompute the arrival/departure % column (same formula pattern, F{row}/G{row}) ---
$ws.Range("H2").Formula = "=F2/G2"
$ws.Range("H3").Formula = "=F3/G3"
$ws.Range("H4").Formula = "=F4/G4"
$ws.Range("H5").Formula = "=F5/G5"
$ws.Range("H6").Formula = "=F6/G6"
$ws.Range("H7").Formula = "=F7/G7"
$ws.Range("H8").Formula = "=F8/G8"
$ws.Range("H9").Formula = "=F9/G9"
$ws.Range("H10").Formula = "=F10/G10"
$ws.Range("H11").Formula = "=F11/G11"
$ws.Range("H12").Formula = "=F12/G12"
$ws.Range("H13").Formula = "=F13/G13"
$ws.Range("H14").Formula = "=F14/G14"

# --- Row 14 label stays "Total" (shared-string slot just gets renumbered on save) ---
$ws.Range("A14").Value = "Total"

# --- New scratch cell K24, same numeric style as column B:F (thousands format) ---
$ws.Range("K24").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'
$ws.Range("K24").Value = ""

# --- New column K width (closest reachable width to the recorded best-fit 12.57) ---
$ws.Columns("K").ColumnWidth = 11.6

# --- Selection moved to H3 (matches the author's last clicked cell) ---
$ws.Range("H3").Select()
